$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume table refresh (GitHub Actions data pull).
# D (Price) and B/C (Coin/Link, rows 45-46 swap) are stored as plain text
# in the source sheet, so force Text number-format before writing the
# Price column to stop Excel from coercing "196.80" -> 196.8, "0.140" -> 0.14, etc.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.836.66"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.497.62"
$ws.Range("E3").Value = "  -1.08%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.96"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.80"
$ws.Range("E6").Value = "  +6.66%  "

$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("E9").Value = "  -1.64%  "

$ws.Range("E10").Value = "  +2.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.22"
$ws.Range("E11").Value = "  +1.41%  "

$ws.Range("E12").Value = "  -1.79%  "

$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.054.25"
$ws.Range("E14").Value = "  -1.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "601.94"
$ws.Range("E15").Value = "  +4.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.932.46"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.12"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.503.43"
$ws.Range("E19").Value = "  -1.93%  "

$ws.Range("E20").Value = "  +0.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.993"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.24"
$ws.Range("E22").Value = "  +5.43%  "

$ws.Range("E23").Value = "  +10.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.05"
$ws.Range("E24").Value = "  +4.90%  "

$ws.Range("E25").Value = "  -1.90%  "

$ws.Range("E26").Value = "  +5.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.96"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("E28").Value = "  +5.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.57"
$ws.Range("E29").Value = "  +4.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.52"
$ws.Range("E30").Value = "  +24.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("E31").Value = "  +3.21%  "

$ws.Range("E32").Value = "  +4.43%  "

$ws.Range("E33").Value = "  +1.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.69"
$ws.Range("E34").Value = "  +0.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.738.79"
$ws.Range("E35").Value = "  +6.00%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0803"
$ws.Range("E37").Value = "  +3.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "510.06"
$ws.Range("E38").Value = "  -3.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  -8.33%  "

$ws.Range("E40").Value = "  -2.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.63"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").Value = "  +0.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  +0.93%  "

$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.83"
$ws.Range("E45").Value = "  -2.25%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.140"
$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("E47").Value = "  -3.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("E49").Value = "  -4.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.93"
$ws.Range("E50").Value = "  -2.00%  "

$ws.Range("E51").Value = "  -0.27%  "

